$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these price cells keep their exact text representation
# (values that look like plain numbers would otherwise be auto-converted)
$textCells = 'D4','D5','D6','D7','D14','D19','D21','D23','D24','D28','D31','D35','D37','D39','D40','D41','D42','D44','D45','D46','D48','D50','D51'
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '56.180.63'
$ws.Range("E2").Value = '  +2.74%  '

# Row 3
$ws.Range("D3").Value = '2.338.04'
$ws.Range("E3").Value = '  +2.59%  '

# Row 4
$ws.Range("D4").Value = '0.995'
$ws.Range("E4").Value = '  -0.61%  '

# Row 5
$ws.Range("D5").Value = '516.12'
$ws.Range("E5").Value = '  +2.35%  '

# Row 6
$ws.Range("D6").Value = '133.18'
$ws.Range("E6").Value = '  +3.48%  '

# Row 7
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("E8").Value = '  +0.91%  '

# Row 9
$ws.Range("D9").Value = '2.334.16'
$ws.Range("E9").Value = '  +1.71%  '

# Row 10
$ws.Range("E10").Value = '  +6.46%  '

# Row 11
$ws.Range("E11").Value = '  +0.16%  '

# Row 12
$ws.Range("E12").Value = '  +7.43%  '

# Row 13
$ws.Range("E13").Value = '  -1.06%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '23.68'
$ws.Range("E14").Value = '  +1.15%  '

# Row 15
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.732.58'
$ws.Range("E15").Value = '  +1.74%  '

# Row 16
$ws.Range("D16").Value = '56.347.56'
$ws.Range("E16").Value = '  +2.93%  '

# Row 17
$ws.Range("E17").Value = '  +2.02%  '

# Row 18
$ws.Range("D18").Value = '2.326.12'
$ws.Range("E18").Value = '  +0.69%  '

# Row 19
$ws.Range("D19").Value = '10.37'
$ws.Range("E19").Value = '  +0.21%  '

# Row 20
$ws.Range("E20").Value = '  +2.48%  '

# Row 21
$ws.Range("D21").Value = '320.29'
$ws.Range("E21").Value = '  +4.35%  '

# Row 22
$ws.Range("E22").Value = '  +2.09%  '

# Row 23
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.18%  '

# Row 24
$ws.Range("D24").Value = '60.35'

# Row 25
$ws.Range("E25").Value = '  +0.53%  '

# Row 26
$ws.Range("E26").Value = '  +5.33%  '

# Row 27
$ws.Range("E27").Value = '  +3.46%  '

# Row 28
$ws.Range("D28").Value = '170.54'
$ws.Range("E28").Value = '  -0.55%  '

# Row 29
$ws.Range("E29").Value = '  +8.64%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0733'
$ws.Range("E30").Value = '  +4.62%  '

# Row 31
$ws.Range("D31").Value = '1.68'
$ws.Range("E31").Value = '  +3.24%  '

# Row 32
$ws.Range("E32").Value = '  +2.48%  '

# Row 33
$ws.Range("E33").Value = '  +1.47%  '

# Row 34
$ws.Range("E34").Value = '  +0.01%  '

# Row 35
$ws.Range("D35").Value = '0.995'
$ws.Range("E35").Value = '  -0.01%  '

# Row 36
$ws.Range("E36").Value = '  +2.29%  '

# Row 37
$ws.Range("D37").Value = '1.24'
$ws.Range("E37").Value = '  +3.12%  '

# Row 38
$ws.Range("E38").Value = '  +4.33%  '

# Row 39
$ws.Range("D39").Value = '1.52'
$ws.Range("E39").Value = '  +7.15%  '

# Row 40
$ws.Range("D40").Value = '37.48'
$ws.Range("E40").Value = '  +2.92%  '

# Row 41
$ws.Range("D41").Value = '0.378'
$ws.Range("E41").Value = '  +0.86%  '

# Row 42
$ws.Range("D42").Value = '137.45'
$ws.Range("E42").Value = '  +8.36%  '

# Row 43
$ws.Range("E43").Value = '  +4.85%  '

# Row 44
$ws.Range("D44").Value = '273.56'
$ws.Range("E44").Value = '  +8.71%  '

# Row 45
$ws.Range("D45").Value = '5.01'
$ws.Range("E45").Value = '  -0.78%  '

# Row 46
$ws.Range("D46").Value = '0.0928'
$ws.Range("E46").Value = '  +3.09%  '

# Row 47
$ws.Range("E47").Value = '  +0.74%  '

# Row 48
$ws.Range("D48").Value = '0.554'
$ws.Range("E48").Value = '  +1.03%  '

# Row 50
$ws.Range("D50").Value = '0.379'
$ws.Range("E50").Value = '  +1.06%  '

# Row 51
$ws.Range("D51").Value = '16.67'
$ws.Range("E51").Value = '  +1.06%  '
